$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Bear Growl / Bear Roar descriptions
$ws.Range("B12").Value = "bear roaring when attacking"
$ws.Range("B11").Value = "bear growling when idling"

# Update Status column (E) for rows whose tracking status changed
# (Made/DUE -> FIXING BUGS, except row 14 which is now fully DONE)
$ws.Range("E5").Value  = "FIXING BUGS"
$ws.Range("E6").Value  = "FIXING BUGS"
$ws.Range("E7").Value  = "FIXING BUGS"
$ws.Range("E11").Value = "FIXING BUGS"
$ws.Range("E12").Value = "FIXING BUGS"
$ws.Range("E14").Value = "DONE"
$ws.Range("E17").Value = "FIXING BUGS"
$ws.Range("E18").Value = "FIXING BUGS"
$ws.Range("E24").Value = "FIXING BUGS"
$ws.Range("E25").Value = "FIXING BUGS"
$ws.Range("E27").Value = "FIXING BUGS"

# Remove the now-unused "DUE" column entirely
$ws.Columns("F").ClearContents()

# Widen the Status column slightly to fit the new values
$ws.Columns("E").ColumnWidth = 12.14

# Restore cursor position as saved by the author
$ws.Range("F25").Select()
